$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 739, pushing existing rows 739-858 down to 741-860
$ws.Range("A739:A740").EntireRow.Insert()

# Row 739 - new "Primera" record for 2023-01-25
$ws.Range("A739").Value = 9
$ws.Range("B739").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C739").Value = "Metropolitana"
$ws.Range("D739").Value = 44951
$ws.Range("D739").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E739").Value = 13
$ws.Range("F739").Value = 100112023
$ws.Range("G739").Value = "Brócoli"
$ws.Range("H739").Value = "Sin especificar"
$ws.Range("I739").Value = "Primera"
$ws.Range("J739").Value = 3400
$ws.Range("K739").Value = 700
$ws.Range("L739").Value = 800
$ws.Range("M739").Value = 750
$ws.Range("N739").Value = "`$/unidad"
$ws.Range("O739").Value = "Región Metropolitana"
$ws.Range("P739").Value = 750
$ws.Range("Q739").Value = 1
$ws.Range("R739").Value = "Hortaliza"

# Row 740 - new "Segunda" record for 2023-01-25
$ws.Range("A740").Value = 9
$ws.Range("B740").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C740").Value = "Metropolitana"
$ws.Range("D740").Value = 44951
$ws.Range("D740").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E740").Value = 13
$ws.Range("F740").Value = 100112023
$ws.Range("G740").Value = "Brócoli"
$ws.Range("H740").Value = "Sin especificar"
$ws.Range("I740").Value = "Segunda"
$ws.Range("J740").Value = 1600
$ws.Range("K740").Value = 600
$ws.Range("L740").Value = 600
$ws.Range("M740").Value = 600
$ws.Range("N740").Value = "`$/unidad"
$ws.Range("O740").Value = "Región Metropolitana"
$ws.Range("P740").Value = 600
$ws.Range("Q740").Value = 1
$ws.Range("R740").Value = "Hortaliza"
